$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 3984.5
$ws.Range("J97").Value = 3984.5
$ws.Range("L97").Value = 11953.5
$ws.Range("N97").Value = -12945.5

$ws.Range("H112").Value = 2391.0588
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 2448.3635
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 7345.0905
$ws.Range("N112").Value = -9561.0905
$ws.Range("M112").Value = -392

$ws.Range("H132").Value = 1704.9231
$ws.Range("I132").Value = 1456.7446
$ws.Range("J132").Value = 4037.8
$ws.Range("K132").Value = 4370.2338
$ws.Range("L132").Value = 12113.4
$ws.Range("M132").Value = -1840.2338
$ws.Range("N132").Value = -17173.4

$ws.Range("H138").Value = 6338.5327
$ws.Range("J138").Value = 6951.662
$ws.Range("L138").Value = 20854.986
$ws.Range("N138").Value = -31134.986

$ws.Range("H141").Value = 5768.2188
$ws.Range("I141").Value = 5722.7666
$ws.Range("K141").Value = 17168.2998
$ws.Range("M141").Value = -11988.2998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 126553.875
$ws.Range("I2").Value = 143918.72
$ws.Range("K2").Value = 143918.72
$ws.Range("M2").Value = -143805.72

$ws.Range("H4").Value = 309.1111
$ws.Range("I4").Value = 131.16667
$ws.Range("K4").Value = 131.16667
$ws.Range("M4").Value = -15.16667000000001

$ws.Range("H32").Value = 5711.5957
$ws.Range("I32").Value = 4305.229
$ws.Range("K32").Value = 4305.229
$ws.Range("M32").Value = -4018.229

$ws.Range("H74").Value = 2009.6034
$ws.Range("J74").Value = 1718.625
$ws.Range("L74").Value = 1718.625
$ws.Range("N74").Value = -3466.625

$ws.Range("H77").Value = 2009.6034
$ws.Range("J77").Value = 1718.625
$ws.Range("L77").Value = 8593.125
$ws.Range("N77").Value = -17329.125

$ws.Range("H116").Value = 126553.875
$ws.Range("I116").Value = 143918.72
$ws.Range("K116").Value = 143918.72
$ws.Range("M116").Value = -141624.72

$ws.Range("H122").Value = 4010.862
$ws.Range("I122").Value = 1924.2307
$ws.Range("K122").Value = 5772.6921
$ws.Range("M122").Value = -3322.6921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 126553.875
$ws.Range("I3").Value = 143918.72
$ws.Range("K3").Value = 143918.72
$ws.Range("M3").Value = -143804.72

$ws.Range("H22").Value = 261.2
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H86").Value = 712523.2
$ws.Range("I86").Value = 898869.8
$ws.Range("K86").Value = 898869.8
$ws.Range("M86").Value = -897746.8

$ws.Range("H89").Value = 712523.2
$ws.Range("I89").Value = 898869.8
$ws.Range("K89").Value = 4494349
$ws.Range("M89").Value = -4488733

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8805.833000000001
$ws.Range("I16").Value = 2333.5454
$ws.Range("K16").Value = 2333.5454
$ws.Range("M16").Value = -2046.5454

$ws.Range("H41").Value = 14465.2
$ws.Range("I41").Value = 3764.7144
$ws.Range("J41").Value = 39433
$ws.Range("K41").Value = 3764.7144
$ws.Range("L41").Value = 39433
$ws.Range("M41").Value = -3336.7144
$ws.Range("N41").Value = -40289

$ws.Range("H50").Value = 24273.533
$ws.Range("I50").Value = 12809.818
$ws.Range("J50").Value = 55798.75
$ws.Range("K50").Value = 12809.818
$ws.Range("L50").Value = 55798.75
$ws.Range("M50").Value = -12184.818
$ws.Range("N50").Value = -57048.75

$ws.Range("H51").Value = 44950
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H58").Value = 2957.3914
$ws.Range("I58").Value = 2977.7856
$ws.Range("J58").Value = 2925.6667
$ws.Range("K58").Value = 2977.7856
$ws.Range("L58").Value = 2925.6667
$ws.Range("M58").Value = -2774.7856
$ws.Range("N58").Value = -3331.6667

$ws.Range("H59").Value = 40680.89
$ws.Range("J59").Value = 40680.89
$ws.Range("L59").Value = 40680.89
$ws.Range("N59").Value = -42970.89

$ws.Range("H61").Value = 44950
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H62").Value = 3335.7144
$ws.Range("I62").Value = 4100
$ws.Range("J62").Value = 3208.3333
$ws.Range("K62").Value = 4100
$ws.Range("L62").Value = 3208.3333
$ws.Range("M62").Value = -3476
$ws.Range("N62").Value = -4456.3333

$ws.Range("H65").Value = 3335.7144
$ws.Range("I65").Value = 4100
$ws.Range("J65").Value = 3208.3333
$ws.Range("K65").Value = 20500
$ws.Range("L65").Value = 16041.6665
$ws.Range("M65").Value = -17380
$ws.Range("N65").Value = -22281.6665

$ws.Range("H86").Value = 25350.842
$ws.Range("I86").Value = 47311.445
$ws.Range("K86").Value = 47311.445
$ws.Range("M86").Value = -46188.445

$ws.Range("H89").Value = 25350.842
$ws.Range("I89").Value = 47311.445
$ws.Range("K89").Value = 236557.225
$ws.Range("M89").Value = -230941.225

$ws.Range("H99").Value = 5375.769
$ws.Range("I99").Value = 4117
$ws.Range("K99").Value = 4117
$ws.Range("M99").Value = -2619

$ws.Range("H105").Value = 987.5833
$ws.Range("I105").Value = 987.5833
$ws.Range("K105").Value = 987.5833
$ws.Range("M105").Value = 759.4167

$ws.Range("H113").Value = 8805.833000000001
$ws.Range("I113").Value = 2333.5454
$ws.Range("K113").Value = 2333.5454
$ws.Range("M113").Value = -163.5454

$ws.Range("H126").Value = 5375.769
$ws.Range("I126").Value = 4117
$ws.Range("K126").Value = 12351
$ws.Range("M126").Value = -9881

$ws.Range("H136").Value = 2957.3914
$ws.Range("I136").Value = 2977.7856
$ws.Range("J136").Value = 2925.6667
$ws.Range("K136").Value = 8933.356800000001
$ws.Range("L136").Value = 8777.000100000001
$ws.Range("M136").Value = -6383.356800000001
$ws.Range("N136").Value = -13877.0001

$ws.Range("H141").Value = 159831.05
$ws.Range("J141").Value = 157964.12
$ws.Range("L141").Value = 157964.12
$ws.Range("N141").Value = -168324.12

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1797.1666
$ws.Range("I23").Value = 2350
$ws.Range("K23").Value = 7050
$ws.Range("M23").Value = -6815

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 42683.793
$ws.Range("I132").Value = 5478.475
$ws.Range("J132").Value = 157161.69
$ws.Range("K132").Value = 16435.425
$ws.Range("L132").Value = 471485.07
$ws.Range("M132").Value = -13905.425
$ws.Range("N132").Value = -476545.07

$ws.Range("H138").Value = 49833.332
$ws.Range("J138").Value = 49833.332
$ws.Range("L138").Value = 49833.332
$ws.Range("N138").Value = -60113.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5044.485
$ws.Range("J132").Value = 7060.7144
$ws.Range("L132").Value = 21182.1432
$ws.Range("N132").Value = -26242.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3533.3333
$ws.Range("I107").Value = 4850
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 14550
$ws.Range("L107").Value = 2700
$ws.Range("M107").Value = -12630
$ws.Range("N107").Value = -6540

$ws.Range("H132").Value = 31624.697
$ws.Range("I132").Value = 1245.6451
$ws.Range("K132").Value = 3736.9353
$ws.Range("M132").Value = -1206.9353
